# Automatic update of files.
#
# 1) Column C ("Förändrad") on every data row (2..81) moves from
#    2023-10-09 (serial 45208) to 2023-10-13 (serial 45212).
# 2) The HYPERLINK() formulas in columns S:Y (present only on rows
#    2,3,4,5 and 33) get their target filenames rewritten to include a
#    descriptive suffix before the extension. Column Y's folder name is
#    also renamed from "tillsynsmail" to "ti,llsynsmail" (typo kept
#    verbatim, matching upstream).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "Förändrad" date column (C) for every data row ---------
$firstRow = 2
$lastRow = 81
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45212
}

# --- 2) Rewrite the document-link formulas in columns S:Y ---------------
# Folder (unchanged unless noted), and the " suffix" inserted right
# before the file extension for each destination column.
$linkCols = @{
    "S" = @{ Folder = "artfynd";         Suffix = " artfynd" }
    "T" = @{ Folder = "kartor";          Suffix = " karta" }
    "U" = @{ Folder = "knärot";          Suffix = " karta knärot" }
    "V" = @{ Folder = "klagomål";        Suffix = " fsc-klagomål" }
    "W" = @{ Folder = "klagomålsmail";   Suffix = " fsc-klagomål mail" }
    "X" = @{ Folder = "tillsyn";         Suffix = " tillsynsbegäran" }
    "Y" = @{ Folder = "ti,llsynsmail";   Suffix = " tillsynsbegäran mail" }
}

# Original extension per column (before the suffix is appended).
$extCols = @{
    "S" = ".xlsx"
    "T" = ".png"
    "U" = ".png"
    "V" = ".docx"
    "W" = ".docx"
    "X" = ".docx"
    "Y" = ".docx"
}

$baseUrl = "https://klasma.github.io/LoggingDetectiveFiles/Logging_0331"

# Rows that carry link formulas, and which columns are populated on each.
$linkRows = @{
    2  = @("S", "T", "U", "V", "W", "X", "Y")
    3  = @("S", "T", "V", "W", "X", "Y")
    4  = @("S", "T", "V", "W", "X", "Y")
    5  = @("S", "T", "V", "W", "X", "Y")
    33 = @("U", "V", "W", "X", "Y")
}

foreach ($row in $linkRows.Keys) {
    $beteckning = $ws.Range("A$row").Value2
    foreach ($col in $linkRows[$row]) {
        $folder = $linkCols[$col].Folder
        $suffix = $linkCols[$col].Suffix
        $ext = $extCols[$col]
        $url = "$baseUrl/$folder/$beteckning$suffix$ext"
        $formula = '=HYPERLINK("' + $url + '", "' + $beteckning + '")'
        $ws.Range("$col$row").Formula = $formula
    }
}
